# Daily attendance processing - 2025-12-03 10:30:26
# Applies the attendance-session refresh to the "Session Analysis Results" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Recorder-email list re-orderings (values only, formatting unchanged)
# ---------------------------------------------------------------------------
$ws.Range("G2").Value  = "majorelle.magdy@med.asu.edu.eg, nourhan.mahmoud@med.asu.edu.eg, shaimaa.ahmed@med.asu.edu.eg, rana.abozaid@med.asu.edu.eg, servinaz@med.asu.edu.eg"
$ws.Range("G18").Value = "aya.hanafy@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg, yasmin.m.senosy@med.asu.edu.eg, shorokmohamed@med.asu.edu.eg"
$ws.Range("G24").Value = "majorelle.magdy@med.asu.edu.eg, nourhan.mahmoud@med.asu.edu.eg, shaimaa.ahmed@med.asu.edu.eg, rana.abozaid@med.asu.edu.eg, servinaz@med.asu.edu.eg"
$ws.Range("G40").Value = "aya.hanafy@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg, yasmin.m.senosy@med.asu.edu.eg, shorokmohamed@med.asu.edu.eg"
$ws.Range("G52").Value = "mariam.noureldin@med.asu.edu.eg, Shimaa.ashraf@med.asu.edu.eg"
$ws.Range("G62").Value = "aya.hanafy@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg, yassmen.ahmed@med.asu.edu.eg, shorokmohamed@med.asu.edu.eg"
$ws.Range("G74").Value = "mariam.noureldin@med.asu.edu.eg, Shimaa.ashraf@med.asu.edu.eg"
$ws.Range("G84").Value = "aya.hanafy@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg, yassmen.ahmed@med.asu.edu.eg, shorokmohamed@med.asu.edu.eg"
$ws.Range("G96").Value = "mariam.noureldin@med.asu.edu.eg, Sara_nabil@med.asu.edu.eg, norhan.mohamed@med.asu.edu.eg"
$ws.Range("G98").Value = "basma.hamed@med.asu.edu.eg, amany.raafat@med.asu.edu.eg, marwa_mustafa@med.asu.edu.eg, yassmina.fattoh@med.asu.edu.eg, Eman.m.abosakaya@med.asu.edu.eg"
$ws.Range("G106").Value = "Monica.Eshak@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg, nardine.alfonse@med.asu.edu.eg"
$ws.Range("G118").Value = "mariam.noureldin@med.asu.edu.eg, Sara_nabil@med.asu.edu.eg, norhan.mohamed@med.asu.edu.eg"
$ws.Range("G120").Value = "basma.hamed@med.asu.edu.eg, amany.raafat@med.asu.edu.eg, marwa_mustafa@med.asu.edu.eg, yassmina.fattoh@med.asu.edu.eg, Eman.m.abosakaya@med.asu.edu.eg"
$ws.Range("G128").Value = "Monica.Eshak@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg, nardine.alfonse@med.asu.edu.eg"
$ws.Range("G134").Value = "majorelle.magdy@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg"
$ws.Range("G150").Value = "wafaa.ebida@med.asu.edu.eg, nardine.alfonse@med.asu.edu.eg, naema.gomaa@med.asu.edu.eg"
$ws.Range("G156").Value = "majorelle.magdy@med.asu.edu.eg, alshimaa.atef@med.asu.edu.egm, Mohammedeltanany@med.asu.edu.eg"

# ---------------------------------------------------------------------------
# 2) "Class Statistics" block (K6:L10) refresh
# ---------------------------------------------------------------------------
$ws.Range("L6").Value  = 25        # Recorded Sessions
$ws.Range("L7").Value  = 13        # Missing Sessions
$ws.Range("L8").Value  = 138       # Pending Sessions
$ws.Range("L9").Value  = "14.2%"   # Coverage %
$ws.Range("L10").Value = "29.5%"   # Average Attendance %

# ---------------------------------------------------------------------------
# 3) Per-group summary table (K14:S22) refresh
# ---------------------------------------------------------------------------
$ws.Range("P15").Value = 2
$ws.Range("Q15").Value = 17

$ws.Range("P16").Value = 3
$ws.Range("Q16").Value = 17

$ws.Range("P17").Value = 1
$ws.Range("Q17").Value = 18

$ws.Range("P18").Value = 1
$ws.Range("Q18").Value = 18

$ws.Range("O22").Value = 4
$ws.Range("P22").Value = 0
$ws.Range("R22").Value = "18.2%"
$ws.Range("S22").Value = "11.9%"

# ---------------------------------------------------------------------------
# 4) Sessions that flipped from "Pending" -> "Not Recorded"
#    (copy the row format from an existing "Not Recorded" row so the same
#    style index / pink fill is reused, then update the status text)
# ---------------------------------------------------------------------------
$notRecordedTemplate = $ws.Range("A7:I7")
$pendingToNotRecordedRows = @(19, 36, 54, 73)
foreach ($r in $pendingToNotRecordedRows) {
    $notRecordedTemplate.Copy()
    $ws.Range("A$r`:I$r").PasteSpecial(-4122)  # xlPasteFormats
    $ws.Range("I$r").Value = "Not Recorded"
}
$excel.CutCopyMode = $false

# ---------------------------------------------------------------------------
# 5) Row 172 flipped from "Not Recorded" -> "Recorded" with real attendance
# ---------------------------------------------------------------------------
$recordedTemplate = $ws.Range("A2:I2")
$recordedTemplate.Copy()
$ws.Range("A172:I172").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Range("G172").Value = "wafaa.ebida@med.asu.edu.eg, nardine.alfonse@med.asu.edu.eg, naema.gomaa@med.asu.edu.eg"
$ws.Range("H172").Value = "6/226"
$ws.Range("I172").Value = "Recorded"
